# Update "Periodo Mora" values in the Estado de Cuenta sheet:
#   - old period 2507 -> 2508 (rows 16-17)
#   - old period 2508 -> 2509 (rows 18-19)
# and center-align the "Periodo Mora" column (E) for those data rows,
# matching the rest of the table's centered columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Roll the overdue periods forward by one month.
$ws.Range("E16").Value = "2508"
$ws.Range("E17").Value = "2508"
$ws.Range("E18").Value = "2509"
$ws.Range("E19").Value = "2509"

# Center the period column like the surrounding cells.
$ws.Range("E16:E19").HorizontalAlignment = -4108   # xlCenter
